# Update the "Montecreto" report sheet with 4 additional days of data
# (rows 252-255), matching the style/format of the preceding row (251).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New data to append: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newRows = @(
    @{ Row = 252; Date = 44326; B = 1; C = 1; D = 109.1703056768559 },
    @{ Row = 253; Date = 44327; B = 0; C = 1; D = 109.1703056768559 },
    @{ Row = 254; Date = 44328; B = 0; C = 1; D = 109.1703056768559 },
    @{ Row = 255; Date = 44329; B = 1; C = 2; D = 218.3406113537118 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the formatting (number format, alignment, border, font, style index)
    # of the last existing data row (251) for column A, so the new date cell
    # keeps the same style as the rest of the column, then overwrite the value.
    $ws.Range("A251").Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteFormats)
    $ws.Range("A$row").Value = $r.Date

    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

$excel.CutCopyMode = 0
